$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương cơ bản")

# --- Fix a mis-placed marker in row 15: the "1" belongs to the SÓC TRĂNG
# column (Q), not the CẦN THƠ column (O) -------------------------------
$ws.Range("O15").Clear()
$ws.Range("Q14").Copy()
$ws.Range("Q15").PasteSpecial(-4122)
$ws.Range("Q15").Value = 1
$excel.CutCopyMode = $false

# --- New column S: "Sinh hoạt tại cơ sở" -------------------------------
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("S1").Value = "Sinh hoạt tại cơ sở"
$excel.CutCopyMode = $false

$ws.Range("S2").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("S16").Value = 1
$ws.Range("S17").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("S23").Value = 0
$ws.Range("S24").Value = 0
$ws.Range("S25").Value = 0
$ws.Range("S26").Value = 0

$ws.Columns.Item(19).ColumnWidth = 16.6

# --- The "Tổng lương cơ bản" column (N) loses its explicit cell style --
$ws.Range("N2:N26").ClearFormats()

# --- Update the remembered selection -----------------------------------
$ws.Range("T22").Select()
